# issue #5: stock data output to json file
#
# Adds a "property_category" column (value "stock") to the 股票 (stock)
# worksheet, right before the existing "date" column, and fixes a stray
# space in one of the company names.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new column at H; this shifts the existing date /
# legislator_name / legislator_id columns from H/I/J to I/J/K and carries
# their cell styles (header bold/border style + body style) along.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# Every stock row is categorized as "stock".
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
$ws.Range("H5").Value = "stock"
$ws.Range("H6").Value = "stock"
$ws.Range("H7").Value = "stock"

# Fix stray space in "為升電裝工業股份有限公 司" -> "為升電裝工業股份有限公司"
$ws.Range("B2").Value = "為升電裝工業股份有限公司"
